# Updates the cryptocurrency price/volume snapshot data in the worksheet
# to match the latest scrape, per the commit's unified diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: values are entered with a leading apostrophe and the cell style is
# reset to "Normal" afterwards. This forces Excel to keep numeric-looking
# strings (e.g. "213.01", "1.00") as literal text -- matching the original
# inline-string cells -- without leaving a custom text number format behind.

# Row 2: D2, E2
$ws.Range("D2").Value = "'27.695.31"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.07%  "
$ws.Range("E2").Style = "Normal"

# Row 3: D3, E3
$ws.Range("D3").Value = "'1.644.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.23%  "
$ws.Range("E3").Style = "Normal"

# Row 4: E4
$ws.Range("E4").Value = "'  +0.04%  "
$ws.Range("E4").Style = "Normal"

# Row 5: D5, E5
$ws.Range("D5").Value = "'213.01"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.90%  "
$ws.Range("E5").Style = "Normal"

# Row 6: E6
$ws.Range("E6").Value = "'  -0.28%  "
$ws.Range("E6").Style = "Normal"

# Row 8: D8, E8
$ws.Range("D8").Value = "'23.31"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.73%  "
$ws.Range("E8").Style = "Normal"

# Row 9: E9
$ws.Range("E9").Value = "'  +1.37%  "
$ws.Range("E9").Style = "Normal"

# Row 10: E10
$ws.Range("E10").Value = "'  +0.67%  "
$ws.Range("E10").Style = "Normal"

# Row 11: D11, E11
$ws.Range("D11").Value = "'0.0894"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.52%  "
$ws.Range("E11").Style = "Normal"

# Row 12: D12, E12
$ws.Range("D12").Value = "'1.876.84"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.20%  "
$ws.Range("E12").Style = "Normal"

# Row 13: D13, E13
$ws.Range("D13").Value = "'1.650.25"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.77%  "
$ws.Range("E13").Style = "Normal"

# Row 14: E14
$ws.Range("E14").Value = "'  +0.23%  "
$ws.Range("E14").Style = "Normal"

# Row 15: D15, E15
$ws.Range("D15").Value = "'0.561"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.96%  "
$ws.Range("E15").Style = "Normal"

# Row 16: D16, E16
$ws.Range("D16").Value = "'64.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.64%  "
$ws.Range("E16").Style = "Normal"

# Row 17: D17, E17
$ws.Range("D17").Value = "'27.660.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.98%  "
$ws.Range("E17").Style = "Normal"

# Row 18: D18, E18
$ws.Range("D18").Value = "'231.03"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.74%  "
$ws.Range("E18").Style = "Normal"

# Row 19: D19, E19
$ws.Range("D19").Value = "'0.0₃0724"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  +0.81%  "
$ws.Range("E19").Style = "Normal"

# Row 20: D20, E20
$ws.Range("D20").Value = "'7.63"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.94%  "
$ws.Range("E20").Style = "Normal"

# Row 22: E22
$ws.Range("E22").Value = "'  -0.71%  "
$ws.Range("E22").Style = "Normal"

# Row 23: D23, E23
$ws.Range("D23").Value = "'10.01"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +7.45%  "
$ws.Range("E23").Style = "Normal"

# Row 24: E24
$ws.Range("E24").Value = "'  -3.01%  "
$ws.Range("E24").Style = "Normal"

# Row 25: D25, E25
$ws.Range("D25").Value = "'149.85"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.47%  "
$ws.Range("E25").Style = "Normal"

# Row 26: E26
$ws.Range("E26").Value = "'  -0.05%  "
$ws.Range("E26").Style = "Normal"

# Row 27: E27
$ws.Range("E27").Value = "'  -1.25%  "
$ws.Range("E27").Style = "Normal"

# Row 28: B28, C28, D28, E28
$ws.Range("B28").Value = "'EthereumClassic"
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = "'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = "'15.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.82%  "
$ws.Range("E28").Style = "Normal"

# Row 29: B29, C29, D29, E29
$ws.Range("B29").Value = "'BinanceUSD"
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").Value = "'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"

# Row 30: E30
$ws.Range("E30").Value = "'  +0.78%  "
$ws.Range("E30").Style = "Normal"

# Row 31: E31
$ws.Range("E31").Value = "'  +0.73%  "
$ws.Range("E31").Style = "Normal"

# Row 32: E32
$ws.Range("E32").Value = "'  +0.69%  "
$ws.Range("E32").Style = "Normal"

# Row 33: D33, E33
$ws.Range("D33").Value = "'1.444.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  +2.51%  "
$ws.Range("E33").Style = "Normal"

# Row 34: E34
$ws.Range("E34").Value = "'  +1.35%  "
$ws.Range("E34").Style = "Normal"

# Row 35: E35
$ws.Range("E35").Value = "'  +1.73%  "
$ws.Range("E35").Style = "Normal"

# Row 36: E36
$ws.Range("E36").Value = "'  -1.21%  "
$ws.Range("E36").Style = "Normal"

# Row 37: D37, E37
$ws.Range("D37").Value = "'0.568"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +1.19%  "
$ws.Range("E37").Style = "Normal"

# Row 38: D38, E38
$ws.Range("D38").Value = "'0.880"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +0.13%  "
$ws.Range("E38").Style = "Normal"

# Row 39: E39
$ws.Range("E39").Value = "'  +0.11%  "
$ws.Range("E39").Style = "Normal"

# Row 40: D40, E40
$ws.Range("D40").Value = "'0.886"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +12.19%  "
$ws.Range("E40").Style = "Normal"

# Row 41: E41
$ws.Range("E41").Value = "'  +0.24%  "
$ws.Range("E41").Style = "Normal"

# Row 42: E42
$ws.Range("E42").Value = "'  +0.11%  "
$ws.Range("E42").Style = "Normal"

# Row 43: E43
$ws.Range("E43").Value = "'  +3.25%  "
$ws.Range("E43").Style = "Normal"

# Row 44: D44, E44
$ws.Range("D44").Value = "'67.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +4.19%  "
$ws.Range("E44").Style = "Normal"

# Row 45: E45
$ws.Range("E45").Value = "'  +0.61%  "
$ws.Range("E45").Style = "Normal"

# Row 46: D46, E46
$ws.Range("D46").Value = "'1.786.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +0.21%  "
$ws.Range("E46").Style = "Normal"

# Row 47: D47, E47
$ws.Range("D47").Value = "'1.73"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +5.77%  "
$ws.Range("E47").Style = "Normal"

# Row 48: D48, E48
$ws.Range("D48").Value = "'85.68"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.72%  "
$ws.Range("E48").Style = "Normal"

# Row 49: E49
$ws.Range("E49").Value = "'  +0.12%  "
$ws.Range("E49").Style = "Normal"

# Row 50: D50, E50
$ws.Range("D50").Value = "'7.78"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.08%  "
$ws.Range("E50").Style = "Normal"

# Row 51: E51
$ws.Range("E51").Value = "'  +0.89%  "
$ws.Range("E51").Style = "Normal"
